# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (E) and "Correspond Handback
# DateTime" (H) values for the 425d36a0-... file row (row 2) on both the
# "zh-cn" and "de-de" status sheets, reflecting a freshly regenerated
# handback report. The a9377a48-... file row (row 3) is unchanged.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 18:37:52"
$wsZhCn.Range("H2").Value = "2016-03-17 18:38:11"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 18:37:55"
$wsDeDe.Range("H2").Value = "2016-03-17 18:38:19"
